$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values
$ws.Range("B2").Value = 4.9914799125796581
$ws.Range("C2").Value = 10.069361275186626
$ws.Range("D2").Value = 12.036486069528083
$ws.Range("E2").Value = 10.663079852511585

# Row 3 data values
$ws.Range("B3").Value = 3.8442887910512433
$ws.Range("C3").Value = 6.0993965164398682
$ws.Range("D3").Value = 14.79003409347435
$ws.Range("E3").Value = 5.6930233603028739

# Update the selected range to reflect the updated selection in the sheet view
$ws.Range("B1:E3").Select()
